# Se procesan de nuevo los datos con las nuevas dimensiones curadas
#
# The "municipio-nombre" column (D) was previously curated as a measure
# (iaest-measure:municipio-nombre) but is now curated as a dimension
# (sdmx-dimension:refArea / dim / URI-Municipio).
#
# The "numero-de-miembros" column (G) was previously curated as a dimension
# (iaest-dimension:numero-de-miembros / dim / skos:Concept, with a mapping
# file) but is now curated as a measure (iaest-measure:numero-de-miembros /
# medida / xsd:int), so its mapping file reference is removed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D ("municipio-nombre"): measure -> dimension
$ws.Range("D2").Value = "sdmx-dimension:refArea"
$ws.Range("D3").Value = "dim"
$ws.Range("D4").Value = "URI-Municipio"

# Column G ("numero-de-miembros"): dimension -> measure
$ws.Range("G2").Value = "iaest-measure:numero-de-miembros"
$ws.Range("G3").Value = "medida"
$ws.Range("G4").Value = "xsd:int"

# The old dimension mapping file for "numero-de-miembros" no longer applies
$ws.Range("G5").Clear()
